$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - Hong Kong, China (HKG)
$ws.Range("H13").Value = 4.7
$ws.Range("I13").Value = 3.4
$ws.Range("J13").Value = 6
$ws.Range("K13").Value = 4.5
$ws.Range("L13").Value = 6.4
$ws.Range("M13").Value = 6
$ws.Range("N13").ClearContents()
$ws.Range("O13").ClearContents()
$ws.Range("P13").Value = 1.4
$ws.Range("Q13").Value = 1.9

# Row 15 - Indonesia (INO)
$ws.Range("H15").Value = 4.4000000000000004
$ws.Range("I15").Value = 5.4
$ws.Range("J15").Value = 4.0999999999999996
$ws.Range("K15").Value = 5.6
$ws.Range("L15").Value = 4.8
$ws.Range("M15").Value = 6.3
$ws.Range("N15").Value = 3.9
$ws.Range("O15").Value = 4.3
$ws.Range("P15").Value = 2.2999999999999998
$ws.Range("Q15").Value = 3.1

# Autofit the newly populated columns so Excel records bestFit custom widths
$ws.Range("H1:Q48").Columns.AutoFit()

# Update selection to match the recorded UI state after the edit
$ws.Range("R15").Select()
